# Update tests for mobile phone page
# Adds a new "Sheet4" worksheet (mobile phone listing data) after Sheet3,
# makes it the active sheet, and tweaks Sheet3's selection state accordingly.

$wb = $excel.ActiveWorkbook

# --- 1. Add Sheet4 as the last sheet (after Sheet3) and make it active ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)

# --- 2. Populate string cells in the exact order that first introduces each
#        unique value, so the shared-string table comes out in the expected
#        order. ---
$ws4.Range("A1").Value = "For Sale"
$ws4.Range("A11").Value = "Wanted"
$ws4.Range("D4").Value = "null"
$ws4.Range("B1").Value = "Apple"
$ws4.Range("B3").Value = "Samsung"
$ws4.Range("B5").Value = "Xiaomi"
$ws4.Range("B7").Value = "Google"
$ws4.Range("B9").Value = "Vivo"
$ws4.Range("C3").Value = "Galaxy S10"
$ws4.Range("C4").Value = "Galaxy A15"
$ws4.Range("C5").Value = "Redmi Note 13"
$ws4.Range("C7").Value = "Pixel 7 Pro"
$ws4.Range("C8").Value = "Pixel 5"
$ws4.Range("C9").Value = "Y93"
$ws4.Range("D1").Value = "New"
$ws4.Range("D2").Value = "Used"
$ws4.Range("C1").Value = "iPhone 12 Pro"
$ws4.Range("C2").Value = "iPhone 11"

# --- 3. Fill in the remaining cells (reusing already-created shared strings
#        and writing plain numeric values). ---
$ws4.Range("E1").Value = 100000
$ws4.Range("F1").Value = 200000

$ws4.Range("A2").Value = "For Sale"
$ws4.Range("B2").Value = "Apple"
$ws4.Range("E2").Value = 90000
$ws4.Range("F2").Value = 150000

$ws4.Range("A3").Value = "For Sale"
$ws4.Range("D3").Value = "New"
$ws4.Range("E3").Value = 50000
$ws4.Range("F3").Value = 200000

$ws4.Range("A4").Value = "For Sale"
$ws4.Range("B4").Value = "Samsung"
$ws4.Range("E4").Value = "null"
$ws4.Range("F4").Value = 130000

$ws4.Range("A5").Value = "For Sale"
$ws4.Range("D5").Value = "New"
$ws4.Range("E5").Value = 30000
$ws4.Range("F5").Value = "null"

$ws4.Range("A6").Value = "For Sale"
$ws4.Range("B6").Value = "Xiaomi"
$ws4.Range("C6").Value = "null"
$ws4.Range("D6").Value = "Used"
$ws4.Range("E6").Value = "null"
$ws4.Range("F6").Value = "null"

$ws4.Range("A7").Value = "For Sale"
$ws4.Range("D7").Value = "New"
$ws4.Range("E7").Value = 100000
$ws4.Range("F7").Value = 120000

$ws4.Range("A8").Value = "For Sale"
$ws4.Range("B8").Value = "Google"
$ws4.Range("D8").Value = "null"
$ws4.Range("E8").Value = "null"
$ws4.Range("F8").Value = "null"

$ws4.Range("A9").Value = "For Sale"
$ws4.Range("D9").Value = "New"
$ws4.Range("E9").Value = 40000
$ws4.Range("F9").Value = 190000

$ws4.Range("A10").Value = "For Sale"
$ws4.Range("B10").Value = "Vivo"
$ws4.Range("C10").Value = "null"
$ws4.Range("D10").Value = "Used"
$ws4.Range("E10").Value = 40000
$ws4.Range("F10").Value = 95000

$ws4.Range("B11").Value = "null"
$ws4.Range("C11").Value = "null"
$ws4.Range("D11").Value = "null"
$ws4.Range("E11").Value = "null"
$ws4.Range("F11").Value = "null"

# --- 4. Give the "Galaxy A15" cell its own font (Calibri) so a new font /
#        cell style gets registered, matching the extra font added to
#        styles.xml. ---
$ws4.Range("C4").Font.Name = "Calibri"

# --- 5. Column C is a bit wider than the rest. ---
$ws4.Range("C1").ColumnWidth = 16.14

# --- 6. Page setup or  ientation for the new sheet. ---
$ws4.PageSetup.Orientation = 1

# --- 7. Selection on the new sheet is Q5, not A1. ---
[void]$ws4.Range("Q5").Select()
